# B1--and-B2-PowerPoint.pptx edit
#  1) Slide 5's table switches from the custom "Table_0" style to the
#     built-in "Medium Style 2 - Accent 1" table style.
#  2) The deck's theme palette is swapped from the "Integral" (Red Violet)
#     colours to the stock "Office Theme" colours.

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 5 -------------------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{8DF8C99E-6677-4C53-B5AF-909D04D12FB0}")
    }
}

# --- 2) Theme colour scheme: Integral/"Red Violet" -> "Office Theme" ----
function ConvertTo-RGBValue([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches the DrawingML clrScheme slot order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)

$colorScheme = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = ConvertTo-RGBValue $officeThemeColors[$i - 1]
}
